$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45190 = 2023-09-21)
# that was bumped to 45192 (2023-09-23) for every data row (2 through 443).
$newDate = 45192

for ($row = 2; $row -le 443; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
